$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '65.692.34'
$ws.Cells.Item(2, 5).Value = '  +1.97%  '

$ws.Cells.Item(3, 4).Value = '3.464.86'
$ws.Cells.Item(3, 5).Value = '  +0.09%  '

$ws.Cells.Item(4, 5).Value = '  -0.29%  '

$ws.Cells.Item(5, 4).Value = '''581.88'
$ws.Cells.Item(5, 5).Value = '  +1.56%  '

$ws.Cells.Item(6, 4).Value = '''168.35'
$ws.Cells.Item(6, 5).Value = '  +5.30%  '

$ws.Cells.Item(7, 5).Value = '  -0.13%  '

$ws.Cells.Item(8, 4).Value = '3.466.43'
$ws.Cells.Item(8, 5).Value = '  +0.01%  '

$ws.Cells.Item(9, 5).Value = '  -1.17%  '

$ws.Cells.Item(10, 4).Value = '''7.27'
$ws.Cells.Item(10, 5).Value = '  +0.69%  '

$ws.Cells.Item(11, 5).Value = '  +1.38%  '

$ws.Cells.Item(12, 5).Value = '  -1.51%  '

$ws.Cells.Item(13, 4).Value = '4.060.73'
$ws.Cells.Item(13, 5).Value = '  -0.16%  '

$ws.Cells.Item(14, 5).Value = '  +0.13%  '

$ws.Cells.Item(15, 4).Value = '''27.52'
$ws.Cells.Item(15, 5).Value = '  -0.14%  '

$ws.Cells.Item(16, 5).Value = '  +0.02%  '

$ws.Cells.Item(17, 4).Value = '65.572.67'
$ws.Cells.Item(17, 5).Value = '  +1.30%  '

$ws.Cells.Item(18, 4).Value = '3.333.81'
$ws.Cells.Item(18, 5).Value = '  -3.65%  '

$ws.Cells.Item(19, 4).Value = '''6.23'
$ws.Cells.Item(19, 5).Value = '  +0.06%  '

$ws.Cells.Item(20, 4).Value = '''13.80'
$ws.Cells.Item(20, 5).Value = '  -0.26%  '

$ws.Cells.Item(21, 4).Value = '''384.56'
$ws.Cells.Item(21, 5).Value = '  +0.87%  '

$ws.Cells.Item(22, 5).Value = '  -0.22%  '

$ws.Cells.Item(23, 5).Value = '  +0.16%  '

$ws.Cells.Item(24, 4).Value = '''71.60'
$ws.Cells.Item(24, 5).Value = '  -1.35%  '

$ws.Cells.Item(25, 5).Value = '  -1.43%  '

$ws.Cells.Item(26, 5).Value = '  +1.08%  '

$ws.Cells.Item(27, 4).Value = '''9.86'
$ws.Cells.Item(27, 5).Value = '  +0.34%  '

$ws.Cells.Item(28, 5).Value = '  +1.46%  '

$ws.Cells.Item(29, 4).Value = '''0.997'
$ws.Cells.Item(29, 5).Value = '  -0.39%  '

$ws.Cells.Item(30, 4).Value = '''6.26'
$ws.Cells.Item(30, 5).Value = '  +2.12%  '

$ws.Cells.Item(31, 5).Value = '  +0.60%  '

$ws.Cells.Item(32, 5).Value = '  +0.92%  '

$ws.Cells.Item(33, 4).Value = '''23.28'
$ws.Cells.Item(33, 5).Value = '  -0.18%  '

$ws.Cells.Item(34, 4).Value = '''7.31'
$ws.Cells.Item(34, 5).Value = '  +3.78%  '

$ws.Cells.Item(36, 5).Value = '  -3.60%  '

$ws.Cells.Item(37, 4).Value = '''160.38'
$ws.Cells.Item(37, 5).Value = '  -0.57%  '

$ws.Cells.Item(38, 4).Value = '''0.894'
$ws.Cells.Item(38, 5).Value = '  +8.44%  '

$ws.Cells.Item(39, 5).Value = '  -0.56%  '

$ws.Cells.Item(40, 2).Value = 'RenderToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(40, 4).Value = '''6.64'
$ws.Cells.Item(40, 5).Value = '  +3.08%  '

$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).Value = '''0.0733'
$ws.Cells.Item(41, 5).Value = '  -1.47%  '

$ws.Cells.Item(42, 4).Value = '''26.18'
$ws.Cells.Item(42, 5).Value = '  -2.84%  '

$ws.Cells.Item(43, 4).Value = '''26.87'
$ws.Cells.Item(43, 5).Value = '  +3.82%  '

$ws.Cells.Item(44, 4).Value = '2.803.60'
$ws.Cells.Item(44, 5).Value = '  -0.99%  '

$ws.Cells.Item(45, 4).Value = '''43.07'
$ws.Cells.Item(45, 5).Value = '  +0.58%  '

$ws.Cells.Item(46, 5).Value = '  -1.14%  '

$ws.Cells.Item(47, 4).Value = '''0.0308'
$ws.Cells.Item(47, 5).Value = '  -0.57%  '

$ws.Cells.Item(48, 4).Value = '''2.46'
$ws.Cells.Item(48, 5).Value = '  +4.15%  '

$ws.Cells.Item(49, 4).Value = '''338.80'
$ws.Cells.Item(49, 5).Value = '  +1.72%  '

$ws.Cells.Item(50, 5).Value = '  +1.33%  '

$ws.Cells.Item(51, 4).Value = '''32.39'
$ws.Cells.Item(51, 5).Value = '  +4.98%  '
